# Generate Report for Handoff
#
# fab95b9e-4b31-485f-9aae-f2b46f2f9f87 finished translation and was
# re-handed-off, so it now sorts to the bottom of the (previously
# alphabetically-ish ordered) rows 7-9 block on every sheet, with a
# refreshed "Ready for handoff" status + new handoff timestamps.
# c4c276b7-... and f1b63e6d-... simply shift up one row each, keeping
# their existing status/timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet: File Name | Path And Name | Extension | Publish URL |
#                 zh-cn | de-de | Latest HO Xliff Generate Date
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A7").Value = "c4c276b7-2d3b-4581-9296-39aaf487959f.md"
$ws.Range("B7").Value = "e2e\c4c276b7-2d3b-4581-9296-39aaf487959f.md"
$ws.Range("E7").Value = "Ready for handoff"
$ws.Range("F7").Value = "Ready for handoff"
$ws.Range("G7").Value = "2016-08-30 08:47:42"

$ws.Range("A8").Value = "f1b63e6d-3787-4383-892d-d67b1bfb93b9.md"
$ws.Range("B8").Value = "e2e\f1b63e6d-3787-4383-892d-d67b1bfb93b9.md"
$ws.Range("G8").Value = "2016-08-30 08:51:59"

$ws.Range("A9").Value = "fab95b9e-4b31-485f-9aae-f2b46f2f9f87.md"
$ws.Range("B9").Value = "e2e\fab95b9e-4b31-485f-9aae-f2b46f2f9f87.md"
$ws.Range("G9").Value = "2016-08-30 09:01:43"

# ---------------------------------------------------------------
# zh-cn sheet: Source File Name (A) | File Extension (B) | Status (C) |
#              ... | Latest Handoff File (G) | Latest Handoff Datetime (H)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A7").Value = "c4c276b7-2d3b-4581-9296-39aaf487959f.md"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("G7").Value = "c4c276b7-2d3b-4581-9296-39aaf487959f.9270e297d57f2cc0ca1d3ad72a094f654c1b207d.zh-cn.xlf"
$ws.Range("H7").Value = "2016-08-30 08:47:37"

$ws.Range("A8").Value = "f1b63e6d-3787-4383-892d-d67b1bfb93b9.md"
$ws.Range("G8").Value = "f1b63e6d-3787-4383-892d-d67b1bfb93b9.ca6854bd753b2321c41d38dcebf24a825ab58506.zh-cn.xlf"
$ws.Range("H8").Value = "2016-08-30 08:51:54"

$ws.Range("A9").Value = "fab95b9e-4b31-485f-9aae-f2b46f2f9f87.md"
$ws.Range("G9").Value = "fab95b9e-4b31-485f-9aae-f2b46f2f9f87.33089da5abdc4ec509960b8f2e1b96998480aff5.zh-cn.xlf"
$ws.Range("H9").Value = "2016-08-30 09:01:28"

# ---------------------------------------------------------------
# de-de sheet: same layout as zh-cn
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A7").Value = "c4c276b7-2d3b-4581-9296-39aaf487959f.md"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("G7").Value = "c4c276b7-2d3b-4581-9296-39aaf487959f.9270e297d57f2cc0ca1d3ad72a094f654c1b207d.de-de.xlf"
$ws.Range("H7").Value = "2016-08-30 08:47:42"

$ws.Range("A8").Value = "f1b63e6d-3787-4383-892d-d67b1bfb93b9.md"
$ws.Range("G8").Value = "f1b63e6d-3787-4383-892d-d67b1bfb93b9.ca6854bd753b2321c41d38dcebf24a825ab58506.de-de.xlf"
$ws.Range("H8").Value = "2016-08-30 08:51:59"

$ws.Range("A9").Value = "fab95b9e-4b31-485f-9aae-f2b46f2f9f87.md"
$ws.Range("G9").Value = "fab95b9e-4b31-485f-9aae-f2b46f2f9f87.33089da5abdc4ec509960b8f2e1b96998480aff5.de-de.xlf"
$ws.Range("H9").Value = "2016-08-30 09:01:43"
